# controlTable_CA.xlsx — "Added the third scenario - Vin Stub"
#
# The control table previously had 4 scenario rows (SELECT/20000101-0,
# SELECT/0-99999999, CHOICE/20000101-0, CHOICE/0-99999999). The sheet is
# reduced down to a single SELECT scenario row whose date range now spans
# the full window (20000101 - 99999999), and the CHOICE rows are removed
# entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 3-5 (the SELECT/0-99999999, CHOICE/20000101-0 and
# CHOICE/0-99999999 scenario rows), leaving only the header (row 1) and the
# first SELECT scenario row (row 2).
$ws.Rows("3:5").Delete() | Out-Null

# Extend the remaining SELECT scenario's expiration date out to the open
# ended sentinel so it now covers the whole effective/expiration window.
$ws.Range("F2").Value = 99999999

# Match the saved selection state.
$ws.Range("E9").Select() | Out-Null
